$wb = $excel.ActiveWorkbook

# Work on the "Users" sheet (sheet3.xml) - add a new row of data.
$ws = $wb.Worksheets.Item("Users")
$ws.Activate() | Out-Null

# New row 5: A5 = "F00474" (plain text), C5 = "074" (text, right-aligned numFmt style like C2:C4)
$ws.Range("A5").Value = "F00474"

# Copy the style from C4 (numFmtId 49, horizontal right) onto C5 before writing the value
# so that it keeps the same cell style (s="4") as the rows above it.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null
$ws.Range("C5").Value = "074"

# Update the selection on the Users sheet to C7 (still within/near the used range)
$ws.Range("C7").Select() | Out-Null
